# Sprint Cadence Diagrams.pptx - "Update Sprint Cadence Diagrams.pptx"
#
# 1. Append a new, blank slide as slide 3 (end of deck).
# 2. Refresh the cached "date last updated" footer text (the
#    datetimeFigureOut field shown in Slide Master view) on the slide
#    master and every slide layout, which PowerPoint re-stamps whenever
#    the deck is edited/saved on a later day.

$p = $ppt.ActivePresentation

# --- 1. Add a new blank slide at the end of the deck ------------------
$newIndex = $p.Slides.Count + 1
$ppLayoutBlank = 12
$newSlide = $p.Slides.Add($newIndex, $ppLayoutBlank)

# --- 2. Bump the cached footer date stamp from 8/20/2020 to 8/21/2020 -
function Update-CachedDateShape {
    param($shapes)

    for ($shapeIdx = 1; $shapeIdx -le $shapes.Count; $shapeIdx++) {
        $shape = $shapes.Item($shapeIdx)
        $text = ""
        try { $text = $shape.TextFrame.TextRange.Text } catch { }
        if ($text -eq "8/20/2020") {
            $shape.TextFrame.TextRange.Text = "8/21/2020"
        }
    }
}

$master = $p.SlideMaster
Update-CachedDateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($layoutIdx = 1; $layoutIdx -le $layouts.Count; $layoutIdx++) {
    Update-CachedDateShape $layouts.Item($layoutIdx).Shapes
}

Write-Output ("Slides: " + $p.Slides.Count)
